$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 263.44446
$ws.Range("I11").Value = 263.44446
$ws.Range("K11").Value = 263.44446
$ws.Range("M11").Value = -123.44446
$ws.Range("H32").Value = 16667077
$ws.Range("I32").Value = 27778128
$ws.Range("K32").Value = 27778128
$ws.Range("M32").Value = -27777802
$ws.Range("H33").Value = 410.1111
$ws.Range("I33").Value = 410.1111
$ws.Range("K33").Value = 410.1111
$ws.Range("M33").Value = -181.1111
$ws.Range("H40").Value = 3948.9697
$ws.Range("J40").Value = 3677.2
$ws.Range("L40").Value = 3677.2
$ws.Range("N40").Value = -4027.2
$ws.Range("H41").Value = 162.9
$ws.Range("I41").Value = 178.5
$ws.Range("K41").Value = 178.5
$ws.Range("M41").Value = 261.5
$ws.Range("H52").Value = 6998
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 6998
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 20994
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -21314
$ws.Range("H58").Value = 1632.75
$ws.Range("I58").Value = 510.5
$ws.Range("K58").Value = 1531.5
$ws.Range("M58").Value = -1381.5
$ws.Range("H70").Value = 4882.6665
$ws.Range("J70").Value = 5459.6
$ws.Range("L70").Value = 16378.8
$ws.Range("N70").Value = -16918.8
$ws.Range("H73").Value = 4882.6665
$ws.Range("J73").Value = 5459.6
$ws.Range("L73").Value = 16378.8
$ws.Range("N73").Value = -18250.8
$ws.Range("H76").Value = 3439.1538
$ws.Range("I76").Value = 3371
$ws.Range("K76").Value = 3371
$ws.Range("M76").Value = -3056
$ws.Range("H79").Value = 3439.1538
$ws.Range("I79").Value = 3371
$ws.Range("K79").Value = 3371
$ws.Range("M79").Value = -2279
$ws.Range("H113").Value = 111113784
$ws.Range("J113").Value = 3999
$ws.Range("L113").Value = 3999
$ws.Range("N113").Value = -10507
$ws.Range("H121").Value = 2132.6667
$ws.Range("J121").Value = 2132.6667
$ws.Range("L121").Value = 6398.000100000001
$ws.Range("N121").Value = -9892.000100000001
$ws.Range("H123").Value = 148920
$ws.Range("J123").Value = 148920
$ws.Range("L123").Value = 148920
$ws.Range("N123").Value = -158720
$ws.Range("H136").Value = 85709
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
$ws.Range("H137").Value = 6464
$ws.Range("I137").Value = 5090
$ws.Range("J137").Value = 8112.8
$ws.Range("K137").Value = 15270
$ws.Range("L137").Value = 24338.4
$ws.Range("M137").Value = -12720
$ws.Range("N137").Value = -29438.4
$ws.Range("H138").Value = 2578.611
$ws.Range("I138").Value = 1149.3334
$ws.Range("J138").Value = 2864.4666
$ws.Range("K138").Value = 3448.0002
$ws.Range("L138").Value = 8593.399800000001
$ws.Range("M138").Value = 1691.9998
$ws.Range("N138").Value = -18873.3998
$ws.Range("H141").Value = 3246.3635
$ws.Range("I141").Value = 3060.5
$ws.Range("K141").Value = 9181.5
$ws.Range("M141").Value = -4001.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1912.6428
$ws.Range("I45").Value = 1773.1666
$ws.Range("K45").Value = 1773.1666
$ws.Range("M45").Value = -1396.1666
$ws.Range("H61").Value = 3046.9
$ws.Range("I61").Value = 1868.8
$ws.Range("J61").Value = 4225
$ws.Range("K61").Value = 1868.8
$ws.Range("L61").Value = 4225
$ws.Range("M61").Value = -1656.8
$ws.Range("N61").Value = -4649
$ws.Range("H74").Value = 234049.33
$ws.Range("I74").Value = 372633.66
$ws.Range("J74").Value = 3075.4443
$ws.Range("K74").Value = 372633.66
$ws.Range("L74").Value = 3075.4443
$ws.Range("M74").Value = -371759.66
$ws.Range("N74").Value = -4823.4443
$ws.Range("H77").Value = 234049.33
$ws.Range("I77").Value = 372633.66
$ws.Range("J77").Value = 3075.4443
$ws.Range("K77").Value = 1863168.3
$ws.Range("L77").Value = 15377.2215
$ws.Range("M77").Value = -1858800.3
$ws.Range("N77").Value = -24113.2215
$ws.Range("H102").Value = 2720.1765
$ws.Range("I102").Value = 2172.6924
$ws.Range("K102").Value = 2172.6924
$ws.Range("M102").Value = -550.6923999999999
$ws.Range("H110").Value = 289.75
$ws.Range("I110").Value = 289.75
$ws.Range("K110").Value = 289.75
$ws.Range("M110").Value = 1755.25
$ws.Range("H122").Value = 3178.6538
$ws.Range("I122").Value = 3093.261
$ws.Range("J122").Value = 3833.3333
$ws.Range("K122").Value = 9279.782999999999
$ws.Range("L122").Value = 11499.9999
$ws.Range("M122").Value = -6829.782999999999
$ws.Range("N122").Value = -16399.9999
$ws.Range("H132").Value = 4900.5264
$ws.Range("I132").Value = 3791
$ws.Range("J132").Value = 6133.3335
$ws.Range("K132").Value = 11373
$ws.Range("L132").Value = 18400.0005
$ws.Range("M132").Value = -8843
$ws.Range("N132").Value = -23460.0005
$ws.Range("H136").Value = 3046.9
$ws.Range("I136").Value = 1868.8
$ws.Range("J136").Value = 4225
$ws.Range("K136").Value = 5606.4
$ws.Range("L136").Value = 12675
$ws.Range("M136").Value = -3056.4
$ws.Range("N136").Value = -17775
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 744.1667
$ws.Range("I22").Value = 740.3333
$ws.Range("J22").Value = 748
$ws.Range("K22").Value = 740.3333
$ws.Range("L22").Value = 748
$ws.Range("M22").Value = -567.3333
$ws.Range("N22").Value = -1094
$ws.Range("H86").Value = 3249.3157
$ws.Range("I86").Value = 3148.8
$ws.Range("J86").Value = 3626.25
$ws.Range("K86").Value = 3148.8
$ws.Range("L86").Value = 3626.25
$ws.Range("M86").Value = -2025.8
$ws.Range("N86").Value = -5872.25
$ws.Range("H89").Value = 3249.3157
$ws.Range("I89").Value = 3148.8
$ws.Range("J89").Value = 3626.25
$ws.Range("K89").Value = 15744
$ws.Range("L89").Value = 18131.25
$ws.Range("M89").Value = -10128
$ws.Range("N89").Value = -29363.25
$ws.Range("H99").Value = 2711
$ws.Range("I99").Value = 1651.3636
$ws.Range("K99").Value = 1651.3636
$ws.Range("M99").Value = -153.3635999999999
$ws.Range("H105").Value = 17335506
$ws.Range("I105").Value = 1112809.5
$ws.Range("K105").Value = 1112809.5
$ws.Range("M105").Value = -1111062.5
$ws.Range("H107").Value = 6993989
$ws.Range("I107").Value = 6993989
$ws.Range("K107").Value = 6993989
$ws.Range("M107").Value = -6992069
$ws.Range("H134").Value = 1461.6571
$ws.Range("I134").Value = 1042.9259
$ws.Range("J134").Value = 2874.875
$ws.Range("K134").Value = 3128.7777
$ws.Range("L134").Value = 8624.625
$ws.Range("M134").Value = -593.7776999999996
$ws.Range("N134").Value = -13694.625
$ws.Range("H138").Value = 59284
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 189.27272
$ws.Range("I7").Value = 120.42857
$ws.Range("K7").Value = 120.42857
$ws.Range("M7").Value = -7.428569999999993
$ws.Range("H16").Value = 1749.3334
$ws.Range("I16").Value = 1749.3334
$ws.Range("K16").Value = 1749.3334
$ws.Range("M16").Value = -1462.3334
$ws.Range("H31").Value = 2454787.5
$ws.Range("J31").Value = 2844934.5
$ws.Range("L31").Value = 2844934.5
$ws.Range("N31").Value = -2845524.5
$ws.Range("H33").Value = 1879.8
$ws.Range("I33").Value = 1879.8
$ws.Range("K33").Value = 1879.8
$ws.Range("M33").Value = -1500.8
$ws.Range("H34").Value = 2454787.5
$ws.Range("J34").Value = 2844934.5
$ws.Range("L34").Value = 2844934.5
$ws.Range("N34").Value = -2845338.5
$ws.Range("H35").Value = 25
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H36").Value = 17500
$ws.Range("J36").Value = 25000
$ws.Range("L36").Value = 25000
$ws.Range("N36").Value = -25776
$ws.Range("H40").Value = 17500
$ws.Range("J40").Value = 25000
$ws.Range("L40").Value = 25000
$ws.Range("N40").Value = -25320
$ws.Range("H50").Value = 59996
$ws.Range("J50").Value = 59996
$ws.Range("L50").Value = 59996
$ws.Range("N50").Value = -61246
$ws.Range("H58").Value = 2446.1365
$ws.Range("I58").Value = 1424.6666
$ws.Range("J58").Value = 3671.9
$ws.Range("K58").Value = 1424.6666
$ws.Range("L58").Value = 3671.9
$ws.Range("M58").Value = -1221.6666
$ws.Range("N58").Value = -4077.9
$ws.Range("H62").Value = 1005
$ws.Range("I62").Value = 1005
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 1005
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -381
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 1005
$ws.Range("I65").Value = 1005
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 5025
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -1905
$ws.Range("N65").ClearContents()
$ws.Range("H113").Value = 1749.3334
$ws.Range("I113").Value = 1749.3334
$ws.Range("K113").Value = 1749.3334
$ws.Range("M113").Value = 420.6666
$ws.Range("H122").Value = 837.5
$ws.Range("I122").Value = 916.6667
$ws.Range("J122").Value = 600
$ws.Range("K122").Value = 2750.0001
$ws.Range("L122").Value = 1800
$ws.Range("M122").Value = -300.0001000000002
$ws.Range("N122").Value = -6700
$ws.Range("H125").Value = 40000
$ws.Range("J125").Value = 40000
$ws.Range("L125").Value = 40000
$ws.Range("N125").Value = -44920
$ws.Range("H132").Value = 4737.442
$ws.Range("I132").Value = 3983.639
$ws.Range("K132").Value = 11950.917
$ws.Range("M132").Value = -9420.917000000001
$ws.Range("H134").Value = 4315.909
$ws.Range("I134").Value = 4500.8057
$ws.Range("K134").Value = 13502.4171
$ws.Range("M134").Value = -10967.4171
$ws.Range("H136").Value = 2446.1365
$ws.Range("I136").Value = 1424.6666
$ws.Range("J136").Value = 3671.9
$ws.Range("K136").Value = 4273.9998
$ws.Range("L136").Value = 11015.7
$ws.Range("M136").Value = -1723.9998
$ws.Range("N136").Value = -16115.7

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 828.9
$ws.Range("J5").Value = 874.25
$ws.Range("L5").Value = 2622.75
$ws.Range("N5").Value = -2846.75
$ws.Range("H32").Value = 7849.3335
$ws.Range("J32").Value = 7849.3335
$ws.Range("L32").Value = 23548.0005
$ws.Range("N32").Value = -24114.0005
$ws.Range("H33").Value = 343.5
$ws.Range("I33").Value = 290
$ws.Range("K33").Value = 1740
$ws.Range("M33").Value = -1457
$ws.Range("H36").Value = 1123.75
$ws.Range("I36").Value = 498.33334
$ws.Range("K36").Value = 1495.00002
$ws.Range("M36").Value = -1326.00002
$ws.Range("H62").Value = 8506.5
$ws.Range("I62").Value = 9999
$ws.Range("J62").Value = 7014
$ws.Range("K62").Value = 29997
$ws.Range("L62").Value = 21042
$ws.Range("M62").Value = -29311
$ws.Range("N62").Value = -22414
$ws.Range("H65").Value = 8506.5
$ws.Range("I65").Value = 9999
$ws.Range("J65").Value = 7014
$ws.Range("K65").Value = 89991
$ws.Range("L65").Value = 63126
$ws.Range("M65").Value = -86559
$ws.Range("N65").Value = -69990
$ws.Range("H74").Value = 29286.715
$ws.Range("J74").Value = 27499
$ws.Range("L74").Value = 82497
$ws.Range("N74").Value = -84619
$ws.Range("H77").Value = 29286.715
$ws.Range("J77").Value = 27499
$ws.Range("L77").Value = 247491
$ws.Range("N77").Value = -258099
$ws.Range("H81").Value = 4248.5
$ws.Range("J81").Value = 2999
$ws.Range("L81").Value = 8997
$ws.Range("N81").Value = -11243
$ws.Range("H84").Value = 4248.5
$ws.Range("J84").Value = 2999
$ws.Range("L84").Value = 26991
$ws.Range("N84").Value = -38223
$ws.Range("H97").Value = 837870
$ws.Range("I97").Value = 1667999.6
$ws.Range("J97").Value = 7740.3335
$ws.Range("K97").Value = 5003998.800000001
$ws.Range("L97").Value = 23221.0005
$ws.Range("M97").Value = -5003502.800000001
$ws.Range("N97").Value = -24213.0005
$ws.Range("H121").Value = 9191823
$ws.Range("I121").Value = 14286531
$ws.Range("K121").Value = 42859593
$ws.Range("M121").Value = -42858283
$ws.Range("H128").Value = 302997
$ws.Range("I128").Value = 302997
$ws.Range("K128").Value = 908991
$ws.Range("M128").Value = -904011
$ws.Range("H135").Value = 828.9
$ws.Range("J135").Value = 874.25
$ws.Range("L135").Value = 7868.25
$ws.Range("N135").Value = -12938.25
$ws.Range("H139").Value = 8677.607
$ws.Range("I139").Value = 12197.3
$ws.Range("K139").Value = 36591.89999999999
$ws.Range("M139").Value = -31451.89999999999
$ws.Range("H140").Value = 2716.8333
$ws.Range("I140").Value = 2716.8333
$ws.Range("K140").Value = 8150.499899999999
$ws.Range("M140").Value = -2970.499899999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 6635.3335
$ws.Range("J7").Value = 6635.3335
$ws.Range("L7").Value = 6635.3335
$ws.Range("N7").Value = -6859.3335
$ws.Range("H8").Value = 6635.3335
$ws.Range("J8").Value = 6635.3335
$ws.Range("L8").Value = 6635.3335
$ws.Range("N8").Value = -6913.3335
$ws.Range("H20").Value = 20005
$ws.Range("I20").Value = 20005
$ws.Range("K20").Value = 20005
$ws.Range("M20").Value = -19760
$ws.Range("H80").Value = 76926980
$ws.Range("I80").Value = 125003410
$ws.Range("J80").Value = 4693
$ws.Range("K80").Value = 125003410
$ws.Range("L80").Value = 4693
$ws.Range("M80").Value = -125002412
$ws.Range("N80").Value = -6689
$ws.Range("H83").Value = 76926980
$ws.Range("I83").Value = 125003410
$ws.Range("J83").Value = 4693
$ws.Range("K83").Value = 625017050
$ws.Range("L83").Value = 23465
$ws.Range("M83").Value = -625012058
$ws.Range("N83").Value = -33449
$ws.Range("H97").Value = 1371.4166
$ws.Range("I97").Value = 1441.5454
$ws.Range("J97").Value = 600
$ws.Range("K97").Value = 1441.5454
$ws.Range("L97").Value = 600
$ws.Range("M97").Value = -945.5454
$ws.Range("N97").Value = -1592
$ws.Range("H122").Value = 5127.75
$ws.Range("I122").Value = 3999.5
$ws.Range("J122").Value = 6256
$ws.Range("K122").Value = 11998.5
$ws.Range("L122").Value = 18768
$ws.Range("M122").Value = -9548.5
$ws.Range("N122").Value = -23668
$ws.Range("H126").Value = 6759.421
$ws.Range("I126").Value = 2068
$ws.Range("J126").Value = 13210.125
$ws.Range("K126").Value = 6204
$ws.Range("L126").Value = 39630.375
$ws.Range("M126").Value = -3734
$ws.Range("N126").Value = -44570.375
$ws.Range("H132").Value = 2348.861
$ws.Range("I132").Value = 2133.6086
$ws.Range("J132").Value = 2729.6924
$ws.Range("K132").Value = 6400.825800000001
$ws.Range("L132").Value = 8189.0772
$ws.Range("M132").Value = -3870.825800000001
$ws.Range("N132").Value = -13249.0772

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 4989
$ws.Range("I3").Value = 4989
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 4989
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -4877
$ws.Range("N3").ClearContents()
$ws.Range("H7").Value = 2310.625
$ws.Range("I7").Value = 2080.8333
$ws.Range("K7").Value = 2080.8333
$ws.Range("M7").Value = -1968.8333
$ws.Range("H15").Value = 4989
$ws.Range("I15").Value = 4989
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 4989
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -4819
$ws.Range("N15").ClearContents()
$ws.Range("H21").Value = 5000
$ws.Range("J21").Value = 5000
$ws.Range("L21").Value = 5000
$ws.Range("N21").Value = -5348
$ws.Range("H22").Value = 119049500
$ws.Range("J22").Value = 333334660
$ws.Range("L22").Value = 333334660
$ws.Range("N22").Value = -333335250
$ws.Range("H24").Value = 30000
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()
$ws.Range("H27").Value = 119049500
$ws.Range("J27").Value = 333334660
$ws.Range("L27").Value = 333334660
$ws.Range("N27").Value = -333334874
$ws.Range("H32").Value = 1300.3334
$ws.Range("I32").Value = 1300.3334
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1300.3334
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -983.3334
$ws.Range("N32").ClearContents()
$ws.Range("H46").Value = 3109.8
$ws.Range("I46").Value = 2099.8
$ws.Range("J46").Value = 4119.8
$ws.Range("K46").Value = 2099.8
$ws.Range("L46").Value = 4119.8
$ws.Range("M46").Value = -1911.8
$ws.Range("N46").Value = -4495.8
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H68").Value = 2999.6667
$ws.Range("I68").Value = 2999.6667
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 2999.6667
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -2250.6667
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 2999.6667
$ws.Range("I71").Value = 2999.6667
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 14998.3335
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -11254.3335
$ws.Range("N71").ClearContents()
$ws.Range("H82").Value = 760.3077
$ws.Range("I82").Value = 760.3077
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 760.3077
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -399.3077
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 760.3077
$ws.Range("I85").Value = 760.3077
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 760.3077
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 487.6923
$ws.Range("N85").ClearContents()
$ws.Range("H100").Value = 5846.533
$ws.Range("I100").Value = 5681.636
$ws.Range("J100").Value = 6300
$ws.Range("K100").Value = 5681.636
$ws.Range("L100").Value = 6300
$ws.Range("M100").Value = -5140.636
$ws.Range("N100").Value = -7382
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H122").Value = 7010.4346
$ws.Range("I122").Value = 5338.3335
$ws.Range("J122").Value = 10145.625
$ws.Range("K122").Value = 16015.0005
$ws.Range("L122").Value = 30436.875
$ws.Range("M122").Value = -13565.0005
$ws.Range("N122").Value = -35336.875
$ws.Range("H126").Value = 2310.625
$ws.Range("I126").Value = 2080.8333
$ws.Range("K126").Value = 6242.499899999999
$ws.Range("M126").Value = -3772.499899999999
$ws.Range("H132").Value = 6705.05
$ws.Range("I132").Value = 4224.5
$ws.Range("K132").Value = 12673.5
$ws.Range("M132").Value = -10143.5
$ws.Range("H136").Value = 6388.3076
$ws.Range("I136").Value = 3406.111
$ws.Range("J136").Value = 13098.25
$ws.Range("K136").Value = 10218.333
$ws.Range("L136").Value = 39294.75
$ws.Range("M136").Value = -7668.332999999999
$ws.Range("N136").Value = -44394.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 305
$ws.Range("I13").Value = 305
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 305
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -165
$ws.Range("N13").ClearContents()
$ws.Range("H39").Value = 48666.332
$ws.Range("J39").Value = 48666.332
$ws.Range("L39").Value = 48666.332
$ws.Range("N39").Value = -49492.332
$ws.Range("H40").Value = 39500
$ws.Range("J40").Value = 39500
$ws.Range("L40").Value = 39500
$ws.Range("N40").Value = -39798
$ws.Range("H42").Value = 47500
$ws.Range("J42").Value = 47500
$ws.Range("L42").Value = 47500
$ws.Range("N42").Value = -48256
$ws.Range("H43").Value = 33699.41
$ws.Range("I43").Value = 14000
$ws.Range("K43").Value = 14000
$ws.Range("M43").Value = -13851
$ws.Range("H92").Value = 34998.75
$ws.Range("J92").Value = 34998.75
$ws.Range("L92").Value = 34998.75
$ws.Range("N92").Value = -39990.75
$ws.Range("H93").Value = 49999
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 49999
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 49999
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = -54991
$ws.Range("H100").Value = 66668370
$ws.Range("I100").Value = 1388.5555
$ws.Range("J100").Value = 166668830
$ws.Range("K100").Value = 2777.111
$ws.Range("L100").Value = 333337660
$ws.Range("M100").Value = -2236.111
$ws.Range("N100").Value = -333338742
$ws.Range("H132").Value = 1719.8462
$ws.Range("I132").Value = 1670
$ws.Range("J132").Value = 1799.6
$ws.Range("K132").Value = 5010
$ws.Range("L132").Value = 5398.799999999999
$ws.Range("M132").Value = -2480
$ws.Range("N132").Value = -10458.8
